$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: field/column headers (now 13 columns, A..M, "tickets" table instead of "clients") ---
$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "glpi_id"
$ws.Range("C1").Value = "title"
$ws.Range("D1").Value = "date"
$ws.Range("E1").Value = "closedate"
$ws.Range("F1").Value = "solvedate"
$ws.Range("G1").Value = "status"
$ws.Range("H1").Value = "description"
$ws.Range("I1").Value = "urgency"
$ws.Range("J1").Value = "impact"
$ws.Range("K1").Value = "priority"
$ws.Range("L1").Value = "type"
$ws.Range("M1").Value = "date_creation"

# --- Row 2: template placeholder tokens bound to the "tickets" table ---
$ws.Range("A2").Value = '${table:tickets.id}'
$ws.Range("B2").Value = '${table:tickets.userGlpiId}'
$ws.Range("C2").Value = '${table:tickets.title}'
$ws.Range("D2").Value = '${table:tickets.date}'
$ws.Range("E2").Value = '${table:tickets.closedate}'
$ws.Range("F2").Value = '${table:tickets.solvedate}'
$ws.Range("G2").Value = '${table:tickets.status}'
$ws.Range("H2").Value = '${table:tickets.description}'
$ws.Range("I2").Value = '${table:tickets.urgency}'
$ws.Range("J2").Value = '${table:tickets.impact}'
$ws.Range("K2").Value = '${table:tickets.priority}'
$ws.Range("L2").Value = '${table:tickets.type}'
$ws.Range("M2").Value = '${table:tickets.dateCreation}'

# Row 2 previously carried a custom (grey) font style (B2:I2); the new layout
# drops the custom font entirely and uses the plain default style everywhere,
# so clear any leftover formatting across the whole used range.
$ws.Range("A1:M2").ClearFormats()

# Match the new page setup (explicit portrait orientation was added).
$ws.PageSetup.Orientation = 1

# Match the new active selection (H3) recorded in the sheet view.
[void]$ws.Range("H3").Select()
